$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target B/C/D/E values for rows 2-51 (row number -> values).
# D and E are forced to text via a leading apostrophe so that
# numeric-looking strings (e.g. "1.00", "240.31") stay text, matching
# the source data (inline strings), instead of Excel auto-converting
# them to numbers.
$rows = @(
    @{ Row=2; B="Bitcoin"; C="https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D="'95.421.05"; E="'  +2.82%  " },
    @{ Row=3; B="Ethereum"; C="https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D="'3.598.92"; E="'  +7.19%  " },
    @{ Row=4; B="TetherUSD"; C="https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D="'1.00"; E="'  +0.11%  " },
    @{ Row=5; B="Solana"; C="https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D="'240.31"; E="'  +3.29%  " },
    @{ Row=6; B="BNB"; C="https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D="'648.91"; E="'  +5.33%  " },
    @{ Row=7; B="XRP"; C="https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D="'1.46"; E="'  +6.65%  " },
    @{ Row=8; B="Dogecoin"; C="https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D="'0.408"; E="'  +5.22%  " },
    @{ Row=9; B="USDC"; C="https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D="'1.00"; E="'  -0.09%  " },
    @{ Row=10; B="Cardano"; C="https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D="'1.00"; E="'  +6.18%  " },
    @{ Row=11; B="LidoStakedEther"; C="https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"; D="'3.600.40"; E="'  +7.29%  " },
    @{ Row=12; B="Avalanche"; C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D="'43.02"; E="'  +0.79%  " },
    @{ Row=13; B="TRON"; C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D="'0.200"; E="'  +2.23%  " },
    @{ Row=14; B="Toncoin"; C="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D="'6.30"; E="'  +1.48%  " },
    @{ Row=15; B="WrappedliquidstakedEther2.0"; C="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D="'4.291.53"; E="'  +7.24%  " },
    @{ Row=16; B="WrappedBTC"; C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D="'95.451.74"; E="'  +2.93%  " },
    @{ Row=17; B="ShibaInu"; C="https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D="'0.0000256"; E="'  +5.00%  " },
    @{ Row=18; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="'3.604.27"; E="'  +7.27%  " },
    @{ Row=19; B="Polkadot"; C="https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D="'7.94"; E="'  -1.51%  " },
    @{ Row=20; B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="'12.51"; E="'  +11.23%  " },
    @{ Row=21; B="Chainlink"; C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D="'18.01"; E="'  +4.01%  " },
    @{ Row=22; B="SuiNetwork"; C="https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"; D="'3.49"; E="'  +4.35%  " },
    @{ Row=23; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="'0.485"; E="'  +13.89%  " },
    @{ Row=24; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="'510.34"; E="'  +3.33%  " },
    @{ Row=25; B="PEPE"; C="https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; D="'0.0000196"; E="'  +7.46%  " },
    @{ Row=26; B="NEARProtocol"; C="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D="'6.65"; E="'  +0.96%  " },
    @{ Row=27; B="Litecoin"; C="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D="'96.57"; E="'  +3.99%  " },
    @{ Row=28; B="Aptos"; C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D="'12.83"; E="'  +8.03%  " },
    @{ Row=29; B="PancakeSwap"; C="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D="'3.12"; E="'  +16.68%  " },
    @{ Row=30; B="InternetComputer(DFINITY)"; C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D="'11.27"; E="'  +1.98%  " },
    @{ Row=31; B="Dai"; C="https://coinranking.com/coin/MoTuySvg7+dai-dai"; D="'1.00"; E="'  -0.06%  " },
    @{ Row=32; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="'0.138"; E="'  +2.52%  " },
    @{ Row=33; B="Binance-PegBSC-USD"; C="https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"; D="'0.992"; E="'  -0.14%  " },
    @{ Row=34; B="Cronos"; C="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D="'0.176"; E="'  +3.01%  " },
    @{ Row=35; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="'31.60"; E="'  +10.99%  " },
    @{ Row=36; B="PolygonEcosystemToken"; C="https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"; D="'0.559"; E="'  +6.88%  " },
    @{ Row=37; B="RenderToken"; C="https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"; D="'8.13"; E="'  +9.26%  " },
    @{ Row=38; B="Bittensor"; C="https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; D="'557.43"; E="'  +0.53%  " },
    @{ Row=39; B="Fetch.AI"; C="https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D="'1.46"; E="'  +6.40%  " },
    @{ Row=40; B="USDe"; C="https://coinranking.com/coin/exbfr2U-0+usde-usde"; D="'1.00"; E="'  -0.12%  " },
    @{ Row=41; B="Kaspa"; C="https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D="'0.150"; E="'  +1.15%  " },
    @{ Row=42; B="ARBITRUM"; C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D="'0.926"; E="'  +5.23%  " },
    @{ Row=43; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="'1.73"; E="'  +1.65%  " },
    @{ Row=44; B="WhiteBITCoin"; C="https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"; D="'23.76"; E="'  +0.38%  " },
    @{ Row=45; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="'5.67"; E="'  +5.66%  " },
    @{ Row=46; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="'0.0418"; E="'  +3.31%  " },
    @{ Row=47; B="Stacks"; C="https://coinranking.com/coin/mMPrMcB7+stacks-stx"; D="'2.25"; E="'  +7.00%  " },
    @{ Row=48; B="OKB"; C="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D="'54.42"; E="'  +3.60%  " },
    @{ Row=49; B="EnergySwap"; C="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D="'32.56"; E="'  +41.75%  " },
    @{ Row=50; B="MantraDAO"; C="https://coinranking.com/coin/cTdD8lD-6+mantradao-om"; D="'3.45"; E="'  -3.67%  " },
    @{ Row=51; B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="'8.10"; E="'  +2.99%  " }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

# The apostrophe-prefixed assignment marks D2:E51 with a "quote prefix"
# cell style; reset those cells back to the Normal style so the cell
# formatting matches the original (unstyled) data cells.
$ws.Range("D2:E51").Style = "Normal"
